$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 939.7778
$ws.Range("I46").Value = 590.2
$ws.Range("J46").Value = 1376.75
$ws.Range("K46").Value = 1770.6
$ws.Range("L46").Value = 4130.25
$ws.Range("M46").Value = -1651.6
$ws.Range("N46").Value = -4368.25

$ws.Range("H60").Value = 939.7778
$ws.Range("I60").Value = 590.2
$ws.Range("J60").Value = 1376.75
$ws.Range("K60").Value = 1770.6
$ws.Range("L60").Value = 4130.25
$ws.Range("M60").Value = -1286.6
$ws.Range("N60").Value = -5098.25

$ws.Range("H62").Value = 1760
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 1450
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 1450
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -2698

$ws.Range("H65").Value = 1760
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 1450
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 7250
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -13490

$ws.Range("H70").Value = 2531.8635
$ws.Range("I70").Value = 2950.25
$ws.Range("J70").Value = 2292.7856
$ws.Range("K70").Value = 8850.75
$ws.Range("L70").Value = 6878.3568
$ws.Range("M70").Value = -8580.75
$ws.Range("N70").Value = -7418.3568

$ws.Range("H73").Value = 2531.8635
$ws.Range("I73").Value = 2950.25
$ws.Range("J73").Value = 2292.7856
$ws.Range("K73").Value = 8850.75
$ws.Range("L73").Value = 6878.3568
$ws.Range("M73").Value = -7914.75
$ws.Range("N73").Value = -8750.356800000001

$ws.Range("H74").Value = 4230.4
$ws.Range("I74").Value = 4080
$ws.Range("J74").Value = 4531.2
$ws.Range("K74").Value = 4080
$ws.Range("L74").Value = 4531.2
$ws.Range("M74").Value = -3144
$ws.Range("N74").Value = -6403.2

$ws.Range("H77").Value = 4230.4
$ws.Range("I77").Value = 4080
$ws.Range("J77").Value = 4531.2
$ws.Range("K77").Value = 20400
$ws.Range("L77").Value = 22656
$ws.Range("M77").Value = -15720
$ws.Range("N77").Value = -32016

$ws.Range("H86").Value = 2083.25
$ws.Range("I86").Value = 2000.5714
$ws.Range("J86").Value = 2199
$ws.Range("K86").Value = 2000.5714
$ws.Range("L86").Value = 2199
$ws.Range("M86").Value = -877.5714
$ws.Range("N86").Value = -4445

$ws.Range("H89").Value = 2083.25
$ws.Range("I89").Value = 2000.5714
$ws.Range("J89").Value = 2199
$ws.Range("K89").Value = 10002.857
$ws.Range("L89").Value = 10995
$ws.Range("M89").Value = -4386.857
$ws.Range("N89").Value = -22227

$ws.Range("H113").Value = 2240.3076
$ws.Range("I113").Value = 2280.4443
$ws.Range("J113").Value = 2150
$ws.Range("K113").Value = 2280.4443
$ws.Range("L113").Value = 2150
$ws.Range("M113").Value = 973.5556999999999
$ws.Range("N113").Value = -8658

$ws.Range("H116").Value = 7350.5
$ws.Range("I116").Value = 8822.143
$ws.Range("K116").Value = 8822.143
$ws.Range("M116").Value = -5380.143

$ws.Range("H126").Value = 23825
$ws.Range("J126").Value = 23825
$ws.Range("L126").Value = 23825
$ws.Range("N126").Value = -33705

$ws.Range("H130").Value = 43472.5
$ws.Range("J130").Value = 43472.5
$ws.Range("L130").Value = 43472.5
$ws.Range("N130").Value = -53512.5

$ws.Range("H141").Value = 1875.9048
$ws.Range("I141").Value = 1904.421
$ws.Range("J141").Value = 1605
$ws.Range("K141").Value = 5713.263
$ws.Range("L141").Value = 4815
$ws.Range("M141").Value = -533.2629999999999
$ws.Range("N141").Value = -15175

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1243.5667
$ws.Range("I74").Value = 1072.28
$ws.Range("J74").Value = 2100
$ws.Range("K74").Value = 1072.28
$ws.Range("L74").Value = 2100
$ws.Range("M74").Value = -198.28
$ws.Range("N74").Value = -3848

$ws.Range("H77").Value = 1243.5667
$ws.Range("I77").Value = 1072.28
$ws.Range("J77").Value = 2100
$ws.Range("K77").Value = 5361.4
$ws.Range("L77").Value = 10500
$ws.Range("M77").Value = -993.3999999999996
$ws.Range("N77").Value = -19236

$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6062394
$ws.Range("I86").Value = 7753631.5
$ws.Range("J86").Value = 2125.0833
$ws.Range("K86").Value = 7753631.5
$ws.Range("L86").Value = 2125.0833
$ws.Range("M86").Value = -7752508.5
$ws.Range("N86").Value = -4371.0833

$ws.Range("H89").Value = 6062394
$ws.Range("I89").Value = 7753631.5
$ws.Range("J89").Value = 2125.0833
$ws.Range("K89").Value = 38768157.5
$ws.Range("L89").Value = 10625.4165
$ws.Range("M89").Value = -38762541.5
$ws.Range("N89").Value = -21857.4165

$ws.Range("H99").Value = 66667772
$ws.Range("I99").Value = 100000720
$ws.Range("J99").Value = 1873.8
$ws.Range("K99").Value = 100000720
$ws.Range("L99").Value = 1873.8
$ws.Range("M99").Value = -99999222
$ws.Range("N99").Value = -4869.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1575.1666
$ws.Range("I16").Value = 977.8
$ws.Range("J16").Value = 2001.8572
$ws.Range("K16").Value = 977.8
$ws.Range("L16").Value = 2001.8572
$ws.Range("M16").Value = -690.8
$ws.Range("N16").Value = -2575.8572

$ws.Range("H22").Value = 306.8
$ws.Range("I22").Value = 200.2
$ws.Range("J22").Value = 520
$ws.Range("K22").Value = 200.2
$ws.Range("L22").Value = 520
$ws.Range("M22").Value = 149.8
$ws.Range("N22").Value = -1220

$ws.Range("H31").Value = 2679.3547
$ws.Range("I31").Value = 1522
$ws.Range("K31").Value = 1522
$ws.Range("M31").Value = -1227

$ws.Range("H34").Value = 2679.3547
$ws.Range("I34").Value = 1522
$ws.Range("K34").Value = 1522
$ws.Range("M34").Value = -1320

$ws.Range("H99").Value = 4637.4736
$ws.Range("I99").Value = 3218.625
$ws.Range("K99").Value = 3218.625
$ws.Range("M99").Value = -1720.625

$ws.Range("H113").Value = 1575.1666
$ws.Range("I113").Value = 977.8
$ws.Range("J113").Value = 2001.8572
$ws.Range("K113").Value = 977.8
$ws.Range("L113").Value = 2001.8572
$ws.Range("M113").Value = 1192.2
$ws.Range("N113").Value = -6341.8572

$ws.Range("H122").Value = 1054.0714
$ws.Range("I122").Value = 855.1429000000001
$ws.Range("J122").Value = 1253
$ws.Range("K122").Value = 2565.4287
$ws.Range("L122").Value = 3759
$ws.Range("M122").Value = -115.4287000000004
$ws.Range("N122").Value = -8659

$ws.Range("H126").Value = 4637.4736
$ws.Range("I126").Value = 3218.625
$ws.Range("K126").Value = 9655.875
$ws.Range("M126").Value = -7185.875

$ws.Range("H138").Value = 38232.855
$ws.Range("J138").Value = 38232.855
$ws.Range("L138").Value = 38232.855
$ws.Range("N138").Value = -48512.855

$ws.Range("H140").Value = 23803.232
$ws.Range("J140").Value = 23803.232
$ws.Range("L140").Value = 23803.232
$ws.Range("N140").Value = -34163.232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 5483.3335
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 6380
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 19140
$ws.Range("N55").Value = -19494

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 629.75
$ws.Range("I97").Value = 658.0909
$ws.Range("K97").Value = 658.0909
$ws.Range("M97").Value = -162.0909

$ws.Range("H102").Value = 1250.6666
$ws.Range("I102").Value = 1180.8
$ws.Range("J102").Value = 1600
$ws.Range("K102").Value = 1180.8
$ws.Range("L102").Value = 1600
$ws.Range("M102").Value = 441.2
$ws.Range("N102").Value = -4844

$ws.Range("H113").Value = 41667700
$ws.Range("I113").Value = 66667300
$ws.Range("J113").Value = 1694.4445
$ws.Range("K113").Value = 66667300
$ws.Range("L113").Value = 1694.4445
$ws.Range("M113").Value = -66665130
$ws.Range("N113").Value = -6034.4445

$ws.Range("H132").Value = 3753.879
$ws.Range("I132").Value = 4808.533
$ws.Range("J132").Value = 2875
$ws.Range("K132").Value = 14425.599
$ws.Range("L132").Value = 8625
$ws.Range("M132").Value = -11895.599
$ws.Range("N132").Value = -13685

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2021.1765
$ws.Range("I7").Value = 1628.4445
$ws.Range("J7").Value = 2463
$ws.Range("K7").Value = 1628.4445
$ws.Range("L7").Value = 2463
$ws.Range("M7").Value = -1516.4445
$ws.Range("N7").Value = -2687

$ws.Range("H126").Value = 2021.1765
$ws.Range("I126").Value = 1628.4445
$ws.Range("J126").Value = 2463
$ws.Range("K126").Value = 4885.333500000001
$ws.Range("L126").Value = 7389
$ws.Range("M126").Value = -2415.333500000001
$ws.Range("N126").Value = -12329

$ws.Range("H134").Value = 40684.285
$ws.Range("J134").Value = 40684.285
$ws.Range("L134").Value = 40684.285
$ws.Range("N134").Value = -50824.285

$ws.Range("H136").Value = 7555.1816
$ws.Range("I136").Value = 7816.7144
$ws.Range("J136").Value = 7097.5
$ws.Range("K136").Value = 23450.1432
$ws.Range("L136").Value = 21292.5
$ws.Range("M136").Value = -20900.1432
$ws.Range("N136").Value = -26392.5
